$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 8502.125
$ws.Range("I38").Value = 1807.8
$ws.Range("J38").Value = 11545
$ws.Range("K38").Value = 5423.4
$ws.Range("L38").Value = 34635
$ws.Range("M38").Value = -5051.4
$ws.Range("N38").Value = -35379

$ws.Range("H52").Value = 1717.25
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30320

$ws.Range("H69").Value = 14993
$ws.Range("J69").Value = 14993
$ws.Range("L69").Value = 44979
$ws.Range("N69").Value = -46727

$ws.Range("H72").Value = 14993
$ws.Range("J72").Value = 14993
$ws.Range("L72").Value = 134937
$ws.Range("N72").Value = -143673

$ws.Range("H98").Value = 2848.5
$ws.Range("I98").Value = 3112.7856
$ws.Range("J98").Value = 998.5
$ws.Range("K98").Value = 3112.7856
$ws.Range("L98").Value = 998.5
$ws.Range("M98").Value = -1614.7856
$ws.Range("N98").Value = -3994.5

$ws.Range("H112").Value = 1514.6177
$ws.Range("J112").Value = 1628.6666
$ws.Range("L112").Value = 4885.9998
$ws.Range("N112").Value = -7101.9998

$ws.Range("H122").Value = 2848.5
$ws.Range("I122").Value = 3112.7856
$ws.Range("J122").Value = 998.5
$ws.Range("K122").Value = 9338.356800000001
$ws.Range("L122").Value = 2995.5
$ws.Range("M122").Value = -6888.356800000001
$ws.Range("N122").Value = -7895.5

$ws.Range("I131").Value = 5000
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 15000
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -9960
$ws.Range("N131").ClearContents()

$ws.Range("H141").Value = 4837.9614
$ws.Range("I141").Value = 4035.8635
$ws.Range("J141").Value = 9249.5
$ws.Range("K141").Value = 12107.5905
$ws.Range("L141").Value = 27748.5
$ws.Range("M141").Value = -6927.5905
$ws.Range("N141").Value = -38108.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2222.5557
$ws.Range("I2").Value = 1200.6154
$ws.Range("J2").Value = 4879.6
$ws.Range("K2").Value = 1200.6154
$ws.Range("L2").Value = 4879.6
$ws.Range("M2").Value = -1087.6154
$ws.Range("N2").Value = -5105.6

$ws.Range("H54").Value = 78949
$ws.Range("J54").Value = 78949
$ws.Range("L54").Value = 78949
$ws.Range("N54").Value = -80487

$ws.Range("H97").Value = 436.88
$ws.Range("I97").Value = 282.8889
$ws.Range("K97").Value = 282.8889
$ws.Range("M97").Value = 213.1111

$ws.Range("H116").Value = 2222.5557
$ws.Range("I116").Value = 1200.6154
$ws.Range("J116").Value = 4879.6
$ws.Range("K116").Value = 1200.6154
$ws.Range("L116").Value = 4879.6
$ws.Range("M116").Value = 1093.3846
$ws.Range("N116").Value = -9467.6

$ws.Range("H122").Value = 4092.389
$ws.Range("I122").Value = 3942.4482
$ws.Range("K122").Value = 11827.3446
$ws.Range("M122").Value = -9377.3446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2222.5557
$ws.Range("I3").Value = 1200.6154
$ws.Range("J3").Value = 4879.6
$ws.Range("K3").Value = 1200.6154
$ws.Range("L3").Value = 4879.6
$ws.Range("M3").Value = -1086.6154
$ws.Range("N3").Value = -5107.6

$ws.Range("H20").Value = 997.37933
$ws.Range("I20").Value = 1051.6957
$ws.Range("J20").Value = 789.1667
$ws.Range("K20").Value = 1051.6957
$ws.Range("L20").Value = 789.1667
$ws.Range("M20").Value = -804.6957
$ws.Range("N20").Value = -1283.1667

$ws.Range("H29").Value = 1697.6
$ws.Range("I29").Value = 1697.6
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1697.6
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1408.6
$ws.Range("N29").ClearContents()

$ws.Range("H94").Value = 2623.2083
$ws.Range("J94").Value = 3083.1
$ws.Range("L94").Value = 3083.1
$ws.Range("N94").Value = -3985.1

$ws.Range("H134").Value = 6910.5
$ws.Range("I134").Value = 6910.5
$ws.Range("K134").Value = 20731.5
$ws.Range("M134").Value = -18196.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5331.5713
$ws.Range("J99").Value = 6076.4
$ws.Range("L99").Value = 6076.4
$ws.Range("N99").Value = -9072.4

$ws.Range("H122").Value = 2020.3125
$ws.Range("I122").Value = 1985.1428
$ws.Range("J122").Value = 2047.6666
$ws.Range("K122").Value = 5955.428400000001
$ws.Range("L122").Value = 6142.9998
$ws.Range("M122").Value = -3505.428400000001
$ws.Range("N122").Value = -11042.9998

$ws.Range("H126").Value = 5331.5713
$ws.Range("J126").Value = 6076.4
$ws.Range("L126").Value = 18229.2
$ws.Range("N126").Value = -23169.2

$ws.Range("H132").Value = 2972.8696
$ws.Range("I132").Value = 2847.55
$ws.Range("K132").Value = 8542.650000000001
$ws.Range("M132").Value = -6012.650000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H57").Value = 259102.5
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H141").Value = 6691.8
$ws.Range("I141").Value = 3383.6
$ws.Range("K141").Value = 10150.8
$ws.Range("M141").Value = -4970.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2182.6333
$ws.Range("I132").Value = 2085.4827
$ws.Range("K132").Value = 6256.4481
$ws.Range("M132").Value = -3726.4481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1230.3529
$ws.Range("I16").Value = 847
$ws.Range("K16").Value = 847
$ws.Range("M16").Value = -677

$ws.Range("H22").Value = 1032.6428
$ws.Range("J22").Value = 1049.16
$ws.Range("L22").Value = 1049.16
$ws.Range("N22").Value = -1639.16

$ws.Range("H27").Value = 1032.6428
$ws.Range("J27").Value = 1049.16
$ws.Range("L27").Value = 1049.16
$ws.Range("N27").Value = -1263.16

$ws.Range("H40").Value = 6406.4287
$ws.Range("I40").Value = 5645.636
$ws.Range("K40").Value = 5645.636
$ws.Range("M40").Value = -5509.636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 20000
$ws.Range("I31").Value = 20000
$ws.Range("K31").Value = 20000
$ws.Range("M31").Value = -19652

$ws.Range("H126").Value = 3438.2666
$ws.Range("I126").Value = 2996.1
$ws.Range("J126").Value = 4322.6
$ws.Range("K126").Value = 8988.299999999999
$ws.Range("L126").Value = 12967.8
$ws.Range("M126").Value = -6518.299999999999
$ws.Range("N126").Value = -17907.8

$ws.Range("H132").Value = 1326.6364
$ws.Range("I132").Value = 1326.3684
$ws.Range("J132").Value = 1328.3334
$ws.Range("K132").Value = 3979.1052
$ws.Range("L132").Value = 3985.0002
$ws.Range("M132").Value = -1449.1052
$ws.Range("N132").Value = -9045.0002
